# The upstream "Artfynd" export re-ordered two observation records that
# used to sit on rows 19 and 20: the record that was on row 20 is now on
# row 19, and vice versa. Swap the per-record field values between the
# two rows (the handful of columns that happen to hold identical data,
# e.g. the site name/municipality, are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose value differs between the two records.
$cols = @("A","B","D","E","F","G","H","I","J","Q","R","Z","AB")

foreach ($col in $cols) {
    $addr19 = $col + "19"
    $addr20 = $col + "20"
    $v19 = $ws.Range($addr19).Value2
    $v20 = $ws.Range($addr20).Value2
    $ws.Range($addr20).Value2 = $v19
    $ws.Range($addr19).Value2 = $v20
}

# Column L held a blank cell on row 20 only; after the swap it belongs to
# row 19 instead. Relocate the (empty) cell rather than assigning "" to
# it, since writing an empty string clears/removes a cell outright.
$ws.Range("L20").Copy($ws.Range("L19"))
$ws.Range("L20").ClearContents()
